# Collapse the three-run "<id>" / "p099v_aN" / "</id>" sequence into a
# single run reading "<id>p099v_N</id>", keeping the formatting (and the
# w:rPr) of the original "<id>" run. There are six such sequences in the
# document, numbered a1..a6 in source order; they become p099v_1..p099v_6.

$d = $word.ActiveDocument
$cursor = 0

for ($n = 1; $n -le 6; $n++) {

    $oldSuffix = "p099v_a$n</id>"
    $newSuffix = "p099v_$n</id>"

    $rng = $d.Range($cursor, $d.Content.End)
    $found = $rng.Find.Execute("<id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find occurrence $n of '<id>'"
    }

    $idEnd = $rng.End

    # Sanity-check that what immediately follows is the run we expect to replace.
    $tailCheck = $d.Range($idEnd, $idEnd + $oldSuffix.Length)
    if ($tailCheck.Text -ne $oldSuffix) {
        throw "Unexpected text after <id> #$n : [$($tailCheck.Text)] (expected [$oldSuffix])"
    }

    # Grow the "<id>" run so it swallows the new id text (merges into this
    # run because it shares the run's formatting), producing one run:
    # "<id>p099v_N</id>" with the original <id> run's rPr.
    $rng.InsertAfter($newSuffix)

    # Now remove the stale "p099v_aN</id>" text that used to follow (it is
    # now shifted right by the length of the text we just inserted).
    $staleStart = $idEnd + $newSuffix.Length
    $staleEnd = $staleStart + $oldSuffix.Length
    $staleRange = $d.Range($staleStart, $staleEnd)
    $staleRange.Delete()

    # Advance the cursor past the just-written "<id>p099v_N</id>" so the
    # next iteration's search starts beyond it.
    $cursor = $staleStart
}
